# Updated symbol list on Fri Jan  6 05:31:23 UTC 2023 with GitHub Actions
# Refresh the crypto price/volume snapshot: several coin rows were re-ranked
# (their Coin/Link/Price/Volume values shifted down a row) and most
# Price/Volume(1h) figures were updated with the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    # Source workbook stores every data cell as text (inlineStr), even
    # numeric-looking Price/Volume figures (so formatting like trailing
    # zeros, "%", and thousands separators survives verbatim). Forcing
    # the NumberFormat to Text before assignment keeps Excel from
    # re-interpreting the string as a number/percentage, and resetting
    # the Style afterwards avoids leaving a stray number-format behind.
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}


# Row 2
Set-TextCell "D2" "257.01"
Set-TextCell "E2" "-0.65%"

# Row 3
Set-TextCell "D3" "27.10"
Set-TextCell "E3" "1.39%"

# Row 4
Set-TextCell "D4" "4.556"
Set-TextCell "E4" "-5.67%"

# Row 5
Set-TextCell "D5" "0.05894"
Set-TextCell "E5" "-1.29%"

# Row 6
Set-TextCell "D6" "6.631"
Set-TextCell "E6" "-0.84%"

# Row 7
Set-TextCell "D7" "0.8541"
Set-TextCell "E7" "-2.44%"

# Row 8
Set-TextCell "D8" "0.9365"
Set-TextCell "E8" "-1.76%"

# Row 9
Set-TextCell "B9" "WazirX"
Set-TextCell "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D9" "0.1385"
Set-TextCell "E9" "-2.34%"

# Row 10
Set-TextCell "B10" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D10" "0.04886"
Set-TextCell "E10" "36.15%"

# Row 11
Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.07072"
Set-TextCell "E11" "-2.10%"

# Row 12
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.03063"
Set-TextCell "E12" "-2.57%"

# Row 13
Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.09111"
Set-TextCell "E13" "-1.37%"

# Row 14
Set-TextCell "B14" "BitForexToken"
Set-TextCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001530"
Set-TextCell "E14" "-0.63%"

# Row 15
Set-TextCell "B15" "One"
Set-TextCell "C15" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D15" "0.0006030"
Set-TextCell "E15" "-0.78%"

# Row 16
Set-TextCell "D16" "0.006019"
Set-TextCell "E16" "-1.07%"

# Row 17
Set-TextCell "E17" "0.18%"

# Row 18
Set-TextCell "D18" "3.180"
Set-TextCell "E18" "-1.37%"

# Row 19
Set-TextCell "E19" "-1.56%"

# Row 20
Set-TextCell "D20" "0.3050"
Set-TextCell "E20" "-2.76%"

# Row 21
Set-TextCell "E21" "-2.79%"

# Row 22
Set-TextCell "D22" "3.918"
Set-TextCell "E22" "11.22%"

# Row 23
Set-TextCell "D23" "0.04268"
Set-TextCell "E23" "0.92%"

# Row 24
Set-TextCell "D24" "0.001221"
Set-TextCell "E24" "0.00%"

# Row 25
Set-TextCell "E25" "-5.03%"

# Row 26
Set-TextCell "E26" "0.05%"

# Row 27
Set-TextCell "E27" "2.09%"

# Row 40
Set-TextCell "D40" "0.03821"
Set-TextCell "E40" "-0.71%"

# Row 41
Set-TextCell "B41" "BKEXToken"
Set-TextCell "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D41" "0.1101"
Set-TextCell "E41" "-0.29%"

# Row 42
Set-TextCell "B42" "KickToken"
Set-TextCell "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D42" "0.003885"
Set-TextCell "E42" "-34.80%"

# Row 43
Set-TextCell "D43" "0.002340"
Set-TextCell "E43" "1.79%"

# Row 44
Set-TextCell "D44" "0.01381"
Set-TextCell "E44" "31.76%"

# Row 45
Set-TextCell "D45" "0.00005378"
Set-TextCell "E45" "-2.02%"

# Row 47
Set-TextCell "D47" "0.06588"

# Row 48
Set-TextCell "D48" "0.2517"
Set-TextCell "E48" "11,735.84%"
